$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'278.22"
$ws.Range("E2").Value = "'6.67%"

$ws.Range("D3").Value = "'27.24"
$ws.Range("E3").Value = "'-0.99%"

$ws.Range("D4").Value = "'4.789"
$ws.Range("E4").Value = "'1.56%"

$ws.Range("D5").Value = "'0.06250"
$ws.Range("E5").Value = "'0.47%"

$ws.Range("D6").Value = "'6.843"
$ws.Range("E6").Value = "'1.63%"

$ws.Range("D7").Value = "'0.8760"
$ws.Range("E7").Value = "'3.09%"

$ws.Range("D8").Value = "'0.9375"
$ws.Range("E8").Value = "'3.08%"

$ws.Range("D9").Value = "'0.1463"
$ws.Range("E9").Value = "'4.60%"

$ws.Range("D10").Value = "'0.05015"
$ws.Range("E10").Value = "'4.13%"

$ws.Range("D11").Value = "'0.07282"
$ws.Range("E11").Value = "'2.80%"

$ws.Range("D12").Value = "'0.03142"
$ws.Range("E12").Value = "'0.59%"

$ws.Range("D13").Value = "'0.09033"
$ws.Range("E13").Value = "'-0.26%"

$ws.Range("D14").Value = "'0.001544"
$ws.Range("E14").Value = "'1.27%"

$ws.Range("D15").Value = "'0.0006252"
$ws.Range("E15").Value = "'1.78%"

$ws.Range("D16").Value = "'0.006071"
$ws.Range("E16").Value = "'0.19%"

$ws.Range("D17").Value = "'3.467"
$ws.Range("E17").Value = "'0.44%"

$ws.Range("E18").Value = "'3.08%"

$ws.Range("E19").Value = "'3.21%"

$ws.Range("D21").Value = "'0.1309"
$ws.Range("E21").Value = "'-0.08%"

$ws.Range("D22").Value = "'3.854"
$ws.Range("E22").Value = "'-5.80%"

$ws.Range("D23").Value = "'0.04320"
$ws.Range("E23").Value = "'1.43%"

$ws.Range("D24").Value = "'0.001177"
$ws.Range("E24").Value = "'-3.54%"

$ws.Range("D25").Value = "'0.004263"
$ws.Range("E25").Value = "'4.41%"

$ws.Range("E26").Value = "'-0.13%"

$ws.Range("D27").Value = "'0.0001612"
$ws.Range("E27").Value = "'-1.67%"

$ws.Range("D40").Value = "'0.04022"
$ws.Range("E40").Value = "'3.03%"

$ws.Range("D41").Value = "'0.006706"
$ws.Range("E41").Value = "'63.07%"

$ws.Range("D42").Value = "'0.1148"
$ws.Range("E42").Value = "'3.42%"

$ws.Range("D43").Value = "'0.002127"
$ws.Range("E43").Value = "'-0.92%"

$ws.Range("D44").Value = "'0.01346"
$ws.Range("E44").Value = "'0.19%"

$ws.Range("D45").Value = "'0.00005124"
$ws.Range("E45").Value = "'-0.12%"

$ws.Range("E46").Value = "'-0.22%"

$ws.Range("D47").Value = "'2.120"
$ws.Range("E47").Value = "'2,940.41%"

$ws.Range("E48").Value = "'-12.19%"

$ws.Range("D49").Value = "'0.00002096"
$ws.Range("E49").Value = "'-0.22%"

$ws.Range("D50").Value = "'0.0001997"
$ws.Range("E50").Value = "'-0.22%"
